$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B12").Value = "quiz2q11"
$ws.Range("C12").Value = "In his video, Roger Peng mentions a concept that allows mixing of text and code. Is this concept called ""Code and Text Mixture""?"
$ws.Range("D12").Value = "Report Yes or No"
